$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking strings are preserved exactly
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '66.393.61'
$ws.Range("E2").Value = '  -1.65%  '
$ws.Range("D3").Value = '2.503.32'
$ws.Range("E3").Value = '  -4.98%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '581.61'
$ws.Range("E5").Value = '  -2.26%  '
$ws.Range("D6").Value = '171.34'
$ws.Range("E6").Value = '  +2.27%  '
$ws.Range("D8").Value = '0.523'
$ws.Range("E8").Value = '  -1.90%  '
$ws.Range("D9").Value = '2.502.26'
$ws.Range("E9").Value = '  -4.99%  '
$ws.Range("E10").Value = '  -1.44%  '
$ws.Range("E11").Value = '  -0.21%  '
$ws.Range("D12").Value = '0.349'
$ws.Range("E12").Value = '  -3.87%  '
$ws.Range("E13").Value = '  -2.19%  '
$ws.Range("D14").Value = '26.59'
$ws.Range("E14").Value = '  -3.94%  '
$ws.Range("D15").Value = '2.957.11'
$ws.Range("E15").Value = '  -5.05%  '
$ws.Range("E16").Value = '  -3.19%  '
$ws.Range("D17").Value = '66.288.01'
$ws.Range("E17").Value = '  -1.49%  '
$ws.Range("D18").Value = '2.540.13'
$ws.Range("E18").Value = '  -3.57%  '
$ws.Range("D19").Value = '11.20'
$ws.Range("E19").Value = '  -6.55%  '
$ws.Range("E20").Value = '  -4.85%  '
$ws.Range("D21").Value = '346.60'
$ws.Range("E21").Value = '  -2.94%  '
$ws.Range("D22").Value = '4.18'
$ws.Range("E22").Value = '  -3.16%  '
$ws.Range("E23").Value = '  -1.46%  '
$ws.Range("D24").Value = '1.96'
$ws.Range("E24").Value = '  +1.18%  '
$ws.Range("E25").Value = '  -0.03%  '
$ws.Range("D26").Value = '69.48'
$ws.Range("E26").Value = '  -0.68%  '
$ws.Range("D27").Value = '9.96'
$ws.Range("E27").Value = '  -3.64%  '
$ws.Range("E28").Value = '  -0.25%  '
$ws.Range("D29").Value = '2.626.88'
$ws.Range("D30").Value = '0.0₃0973'
$ws.Range("E30").Value = '  -3.50%  '
$ws.Range("D31").Value = '525.95'
$ws.Range("E31").Value = '  -4.02%  '
$ws.Range("D32").Value = '8.06'
$ws.Range("E32").Value = '  +1.45%  '
$ws.Range("E33").Value = '  -2.95%  '
$ws.Range("E34").Value = '  -3.28%  '
$ws.Range("E35").Value = '  -4.50%  '
$ws.Range("D36").Value = '0.999'
$ws.Range("E36").Value = '  -0.03%  '
$ws.Range("E37").Value = '  -3.16%  '
$ws.Range("D38").Value = '156.24'
$ws.Range("D39").Value = '18.57'
$ws.Range("E39").Value = '  -2.22%  '
$ws.Range("E40").Value = '  +0.30%  '
$ws.Range("E41").Value = '  -3.41%  '
$ws.Range("E42").Value = '  -2.07%  '
$ws.Range("E43").Value = '  -2.97%  '
$ws.Range("E44").Value = '  -0.04%  '
$ws.Range("E45").Value = '  +2.94%  '
$ws.Range("D46").Value = '39.45'
$ws.Range("E46").Value = '  -1.48%  '
$ws.Range("D47").Value = '148.05'
$ws.Range("E47").Value = '  -3.12%  '
$ws.Range("D48").Value = '0.555'
$ws.Range("E48").Value = '  -4.34%  '
$ws.Range("E49").Value = '  -3.82%  '
$ws.Range("D50").Value = '0.0₆0271'
$ws.Range("E50").Value = '  -9.88%  '
$ws.Range("D51").Value = '1.71'
$ws.Range("E51").Value = '  +0.64%  '
